$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-12 Wednesday" "2025-11-13 Thursday"

Replace-Text "466×7=3262" "321×2=642"
Replace-Text "468×6=2808" "511×5=2555"
Replace-Text "475×9=4275" "340×5=1700"
Replace-Text "104×6=624" "938×3=2814"
Replace-Text "882×2=1764" "944×8=7552"
Replace-Text "804×2=1608" "782×6=4692"
Replace-Text "668×9=6012" "899×8=7192"
Replace-Text "166×2=332" "966×7=6762"
Replace-Text "455×8=3640" "907×5=4535"
Replace-Text "415×6=2490" "655×5=3275"
Replace-Text "135×7=945" "537×7=3759"
Replace-Text "283×4=1132" "466×4=1864"
Replace-Text "143×6=858" "426×5=2130"
Replace-Text "412×6=2472" "401×9=3609"
Replace-Text "925×5=4625" "201×9=1809"
Replace-Text "999×8=7992" "559×7=3913"
Replace-Text "783×6=4698" "940×8=7520"
Replace-Text "219×3=657" "322×5=1610"
Replace-Text "980×4=3920" "623×9=5607"
Replace-Text "695×7=4865" "618×7=4326"
Replace-Text "734×7=5138" "296×9=2664"
Replace-Text "425×6=2550" "885×7=6195"
Replace-Text "911×4=3644" "114×7=798"
Replace-Text "977×4=3908" "965×8=7720"
Replace-Text "790×5=3950" "277×9=2493"
